$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("gof")
$ws.Range("D2").Value = 15574
$ws.Range("F2").Value = 15612
$ws.Range("G2").Value = 15719
$ws.Range("D3").Value = 15530
$ws.Range("F3").Value = 15628
$ws.Range("G3").Value = 15904

$ws = $wb.Worksheets.Item("facets")
$ws.Range("B2").Value = 850
$ws.Range("B3").Value = 829
$ws.Range("B4").Value = 421

$ws = $wb.Worksheets.Item("Estimates 1-2")
$ws.Range("B2").Value = 0.004
$ws.Range("D2").Value = 0.004
$ws.Range("E2").Value = 0.002
$ws.Range("B3").Value = 0.158
$ws.Range("C3").Value = 0.084
$ws.Range("D3").Value = 0.153
$ws.Range("E3").Value = 3.538
$ws.Range("B4").Value = 0.075
$ws.Range("C4").Value = 0.111
$ws.Range("D4").Value = 0.073
$ws.Range("E4").Value = 0.457
$ws.Range("B5").Value = 0.169
$ws.Range("D5").Value = 0.163
$ws.Range("E5").Value = 1.69
$ws.Range("B6").Value = -0.207
$ws.Range("D6").Value = -0.2
$ws.Range("E6").Value = 2.657
$ws.Range("B7").Value = 0.002
$ws.Range("C7").Value = 0.125
$ws.Range("D7").Value = 0.002
$ws.Range("E7").Value = 0
$ws.Range("B8").Value = 0.02
$ws.Range("C8").Value = 0.123
$ws.Range("D8").Value = 0.019
$ws.Range("E8").Value = 0.026
$ws.Range("B9").Value = 0.066
$ws.Range("C9").Value = 0.124
$ws.Range("D9").Value = 0.064
$ws.Range("E9").Value = 0.283
$ws.Range("B10").Value = 0.022
$ws.Range("C10").Value = 0.123
$ws.Range("D10").Value = 0.021
$ws.Range("E10").Value = 0.032
$ws.Range("B11").Value = -0.138
$ws.Range("C11").Value = 0.123
$ws.Range("D11").Value = -0.134
$ws.Range("E11").Value = 1.259
$ws.Range("B12").Value = 0.305
$ws.Range("C12").Value = 0.127
$ws.Range("D12").Value = 0.295
$ws.Range("E12").Value = 5.768
$ws.Range("B13").Value = -0.083
$ws.Range("C13").Value = 0.133
$ws.Range("D13").Value = -0.08
$ws.Range("E13").Value = 0.389
$ws.Range("B14").Value = -0.266
$ws.Range("C14").Value = 0.123
$ws.Range("E14").Value = 4.677
$ws.Range("D14").Value = -0.257
$ws.Range("B15").Value = 0.227
$ws.Range("C15").Value = 0.122
$ws.Range("D15").Value = 0.22
$ws.Range("E15").Value = 3.462
$ws.Range("B16").Value = 0.055
$ws.Range("C16").Value = 0.129
$ws.Range("D16").Value = 0.053
$ws.Range("E16").Value = 0.182
$ws.Range("B17").Value = 0.409
$ws.Range("C17").Value = 0.468
$ws.Range("D17").Value = 0.396
$ws.Range("E17").Value = 0.764

$ws = $wb.Worksheets.Item("Estimates 1-3")
$ws.Range("B2").Value = 0.005
$ws.Range("C2").Value = 0.099
$ws.Range("D2").Value = 0.005
$ws.Range("E2").Value = 0.003
$ws.Range("B3").Value = 0.17
$ws.Range("C3").Value = 0.082
$ws.Range("D3").Value = 0.164
$ws.Range("E3").Value = 4.298
$ws.Range("B4").Value = -0.129
$ws.Range("C4").Value = 0.103
$ws.Range("D4").Value = -0.125
$ws.Range("E4").Value = 1.569
$ws.Range("B5").Value = -0.29
$ws.Range("C5").Value = 0.127
$ws.Range("D5").Value = -0.281
$ws.Range("E5").Value = 5.214
$ws.Range("B6").Value = -0.351
$ws.Range("C6").Value = 0.126
$ws.Range("D6").Value = -0.34
$ws.Range("E6").Value = 7.76
$ws.Range("B7").Value = 0.426
$ws.Range("C7").Value = 0.122
$ws.Range("D7").Value = 0.412
$ws.Range("E7").Value = 12.193
$ws.Range("B8").Value = -0.101
$ws.Range("C8").Value = 0.12
$ws.Range("D8").Value = -0.098
$ws.Range("E8").Value = 0.708
$ws.Range("B9").Value = 0.658
$ws.Range("C9").Value = 0.12
$ws.Range("D9").Value = 0.637
$ws.Range("E9").Value = 30.067
$ws.Range("B10").Value = 0.414
$ws.Range("C10").Value = 0.118
$ws.Range("D10").Value = 0.401
$ws.Range("E10").Value = 12.309
$ws.Range("B11").Value = -0.365
$ws.Range("C11").Value = 0.119
$ws.Range("D11").Value = -0.353
$ws.Range("E11").Value = 9.408
$ws.Range("B12").Value = 0.669
$ws.Range("C12").Value = 0.121
$ws.Range("D12").Value = 0.647
$ws.Range("E12").Value = 30.569
$ws.Range("B13").Value = 0.244
$ws.Range("C13").Value = 0.125
$ws.Range("D13").Value = 0.236
$ws.Range("E13").Value = 3.81
$ws.Range("B14").Value = -0.68
$ws.Range("C14").Value = 0.116
$ws.Range("D14").Value = -0.658
$ws.Range("E14").Value = 34.364
$ws.Range("I14").Value = 0.999
$ws.Range("B15").Value = 0.186
$ws.Range("C15").Value = 0.116
$ws.Range("D15").Value = 0.18
$ws.Range("E15").Value = 2.571
$ws.Range("B16").Value = -0.159
$ws.Range("C16").Value = 0.12
$ws.Range("D16").Value = -0.154
$ws.Range("E16").Value = 1.756
$ws.Range("B17").Value = 0.698
$ws.Range("C17").Value = 0.45
$ws.Range("D17").Value = 0.675
$ws.Range("E17").Value = 2.406

$ws = $wb.Worksheets.Item("Estimates 2-3")
$ws.Range("B2").Value = 0.002
$ws.Range("D2").Value = 0.002
$ws.Range("E2").Value = 0
$ws.Range("B3").Value = 0.012
$ws.Range("D3").Value = 0.012
$ws.Range("E3").Value = 0.019
$ws.Range("B4").Value = -0.205
$ws.Range("C4").Value = 0.119
$ws.Range("D4").Value = -0.198
$ws.Range("E4").Value = 2.968
$ws.Range("B5").Value = -0.459
$ws.Range("C5").Value = 0.133
$ws.Range("D5").Value = -0.444
$ws.Range("E5").Value = 11.91
$ws.Range("B6").Value = -0.144
$ws.Range("D6").Value = -0.139
$ws.Range("E6").Value = 1.246
$ws.Range("B7").Value = 0.425
$ws.Range("C7").Value = 0.128
$ws.Range("D7").Value = 0.411
$ws.Range("E7").Value = 11.024
$ws.Range("B8").Value = -0.121
$ws.Range("C8").Value = 0.126
$ws.Range("D8").Value = -0.117
$ws.Range("E8").Value = 0.922
$ws.Range("B9").Value = 0.591
$ws.Range("C9").Value = 0.128
$ws.Range("D9").Value = 0.572
$ws.Range("E9").Value = 21.318
$ws.Range("B10").Value = 0.392
$ws.Range("C10").Value = 0.128
$ws.Range("D10").Value = 0.379
$ws.Range("E10").Value = 9.379
$ws.Range("B11").Value = -0.227
$ws.Range("D11").Value = -0.22
$ws.Range("E11").Value = 3.145
$ws.Range("B12").Value = 0.364
$ws.Range("C12").Value = 0.133
$ws.Range("D12").Value = 0.352
$ws.Range("E12").Value = 7.49
$ws.Range("B13").Value = 0.327
$ws.Range("C13").Value = 0.14
$ws.Range("D13").Value = 0.316
$ws.Range("E13").Value = 5.456
$ws.Range("B14").Value = -0.413
$ws.Range("C14").Value = 0.13
$ws.Range("D14").Value = -0.4
$ws.Range("E14").Value = 10.093
$ws.Range("B15").Value = -0.042
$ws.Range("C15").Value = 0.129
$ws.Range("D15").Value = -0.041
$ws.Range("E15").Value = 0.106
$ws.Range("B16").Value = -0.213
$ws.Range("C16").Value = 0.138
$ws.Range("D16").Value = -0.206
$ws.Range("E16").Value = 2.382
$ws.Range("B17").Value = 0.289
$ws.Range("C17").Value = 0.486
$ws.Range("D17").Value = 0.28
$ws.Range("E17").Value = 0.354

$ws = $wb.Worksheets.Item("Main effect 1-2")
$ws.Range("B2").Value = 0.727
$ws.Range("C2").Value = 0.704
$ws.Range("B3").Value = 0.707
$ws.Range("C3").Value = 0.686

$ws = $wb.Worksheets.Item("Main effect 1-3")
$ws.Range("B2").Value = 0.305
$ws.Range("C2").Value = 0.295
$ws.Range("B3").Value = 0.298
$ws.Range("C3").Value = 0.289

$ws = $wb.Worksheets.Item("Main effect 2-3")
$ws.Range("B2").Value = -0.422
$ws.Range("C2").Value = -0.408
$ws.Range("B3").Value = -0.409
$ws.Range("C3").Value = -0.397
